$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 174 updates
$ws.Range("B174").Value = 299466
$ws.Range("E174").Value = 53461
$ws.Range("F174").Value = 31923
$ws.Range("G174").Value = 40750
$ws.Range("H174").Value = 32530
$ws.Range("I174").Value = 117480
$ws.Range("J174").Value = 62522
$ws.Range("P174").Value = 6651
$ws.Range("Q174").Value = 28433
$ws.Range("R174").Value = 236944
$ws.Range("U174").Value = 42957
$ws.Range("V174").Value = 27774
$ws.Range("W174").Value = 33642
$ws.Range("X174").Value = 25880
$ws.Range("Y174").Value = 89048

# Row 175 updates
$ws.Range("B175").Value = 286057
$ws.Range("E175").Value = 41735
$ws.Range("F175").Value = 30558
$ws.Range("G175").Value = 41906
$ws.Range("H175").Value = 32644
$ws.Range("I175").Value = 117559
$ws.Range("J175").Value = 59078
$ws.Range("P175").Value = 6351
$ws.Range("Q175").Value = 27930
$ws.Range("R175").Value = 226979
$ws.Range("U175").Value = 33128
$ws.Range("V175").Value = 26408
$ws.Range("W175").Value = 35018
$ws.Range("X175").Value = 26293
$ws.Range("Y175").Value = 89629

# Row 176 updates
$ws.Range("B176").Value = 278789
$ws.Range("D176").Value = 12270
$ws.Range("E176").Value = 41157
$ws.Range("F176").Value = 29084
$ws.Range("G176").Value = 41305
$ws.Range("H176").Value = 33955
$ws.Range("I176").Value = 118171
$ws.Range("J176").Value = 56397
$ws.Range("P176").Value = 6905
$ws.Range("Q176").Value = 27456
$ws.Range("R176").Value = 222392
$ws.Range("T176").Value = 8808
$ws.Range("U176").Value = 33354
$ws.Range("V176").Value = 25479
$ws.Range("W176").Value = 34822
$ws.Range("X176").Value = 27050
$ws.Range("Y176").Value = 90715

# Row 177: previously only had A177, B177, J177 populated.
# Now fill the full row (B through Y) consistent with the other rows.
$ws.Range("B177").Value = 283152
$ws.Range("C177").Value = 2119
$ws.Range("D177").Value = 16986
$ws.Range("E177").Value = 42215
$ws.Range("F177").Value = 27974
$ws.Range("G177").Value = 41757
$ws.Range("H177").Value = 32924
$ws.Range("I177").Value = 119177
$ws.Range("J177").Value = 56850
$ws.Range("K177").Value = 264
$ws.Range("L177").Value = 3561
$ws.Range("M177").Value = 8979
$ws.Range("N177").Value = 3595
$ws.Range("O177").Value = 6496
$ws.Range("P177").Value = 6584
$ws.Range("Q177").Value = 27371
$ws.Range("R177").Value = 226302
$ws.Range("S177").Value = 1855
$ws.Range("T177").Value = 13425
$ws.Range("U177").Value = 33236
$ws.Range("V177").Value = 24379
$ws.Range("W177").Value = 35261
$ws.Range("X177").Value = 26340
$ws.Range("Y177").Value = 91806
